# Add Email Address (P_Email / column BD) for SubmitPPIPayment rows 2-5
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SaleUpgrade")

$ws.Range("BD2:BD5").Value = "iahmed@govolution.com"

# Widen column BD so the new email text is fully visible
$ws.Range("BD1").EntireColumn.ColumnWidth = 27.15234375
